$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ALC")
$ws1.Range("H40").Value = 3160.5386
$ws1.Range("J40").Value = 3245
$ws1.Range("L40").Value = 3245
$ws1.Range("N40").Value = -3595
$ws1.Range("H58").Value = 325
$ws1.Range("I58").Value = 275
$ws1.Range("J58").Value = 550
$ws1.Range("K58").Value = 825
$ws1.Range("L58").Value = 1650
$ws1.Range("M58").Value = -675
$ws1.Range("N58").Value = -1950
$ws1.Range("H76").Value = 7671367
$ws1.Range("J76").Value = 10547029
$ws1.Range("L76").Value = 10547029
$ws1.Range("N76").Value = -10547659
$ws1.Range("H79").Value = 7671367
$ws1.Range("J79").Value = 10547029
$ws1.Range("L79").Value = 10547029
$ws1.Range("N79").Value = -10549213
$ws1.Range("H82").Value = 5000
$ws1.Range("J82").Value = 0
$ws1.Range("L82").Value = 0
$ws1.Range("N82").ClearContents()
$ws1.Range("H85").Value = 5000
$ws1.Range("J85").Value = 0
$ws1.Range("L85").Value = 0
$ws1.Range("N85").ClearContents()
$ws1.Range("H111").Value = 24840.21
$ws1.Range("I111").Value = 1971.4445
$ws1.Range("J111").Value = 45422.1
$ws1.Range("K111").Value = 5914.333500000001
$ws1.Range("L111").Value = 136266.3
$ws1.Range("M111").Value = -2847.333500000001
$ws1.Range("N111").Value = -142400.3
$ws1.Range("H113").Value = 62503780
$ws1.Range("I113").Value = 111114610
$ws1.Range("J113").Value = 4141.2856
$ws1.Range("K113").Value = 111114610
$ws1.Range("L113").Value = 4141.2856
$ws1.Range("M113").Value = -111111356
$ws1.Range("N113").Value = -10649.2856
$ws1.Range("H132").Value = 3354.423
$ws1.Range("J132").Value = 5339.6
$ws1.Range("L132").Value = 16018.8
$ws1.Range("N132").Value = -21078.8
$ws1.Range("H135").Value = 1025.8182
$ws1.Range("I135").Value = 798.7646999999999
$ws1.Range("J135").Value = 1797.8
$ws1.Range("K135").Value = 7188.882299999999
$ws1.Range("L135").Value = 16180.2
$ws1.Range("M135").Value = -4653.882299999999
$ws1.Range("N135").Value = -21250.2
$ws1.Range("H137").Value = 2528.4285
$ws1.Range("I137").Value = 2402.2083
$ws1.Range("K137").Value = 7206.624899999999
$ws1.Range("M137").Value = -4656.624899999999
$ws1.Range("H138").Value = 1844.683
$ws1.Range("I138").Value = 1558.55
$ws1.Range("J138").Value = 2117.1904
$ws1.Range("K138").Value = 4675.65
$ws1.Range("L138").Value = 6351.5712
$ws1.Range("M138").Value = 464.3500000000004
$ws1.Range("N138").Value = -16631.5712

$ws2 = $wb.Worksheets.Item("ARM")
$ws2.Range("H5").Value = 8145.8
$ws2.Range("I5").Value = 243
$ws2.Range("K5").Value = 243
$ws2.Range("M5").Value = -131
$ws2.Range("H45").Value = 1129.9166
$ws2.Range("I45").Value = 997.6
$ws2.Range("K45").Value = 997.6
$ws2.Range("M45").Value = -620.6
$ws2.Range("H61").Value = 3585.0588
$ws2.Range("I61").Value = 1999.091
$ws2.Range("K61").Value = 1999.091
$ws2.Range("M61").Value = -1787.091
$ws2.Range("H102").Value = 2075594.8
$ws2.Range("I102").Value = 2565857.5
$ws2.Range("K102").Value = 2565857.5
$ws2.Range("M102").Value = -2564235.5
$ws2.Range("H136").Value = 3585.0588
$ws2.Range("I136").Value = 1999.091
$ws2.Range("K136").Value = 5997.272999999999
$ws2.Range("M136").Value = -3447.272999999999

$ws3 = $wb.Worksheets.Item("BSM")
$ws3.Range("H4").Value = 8145.8
$ws3.Range("I4").Value = 243
$ws3.Range("K4").Value = 243
$ws3.Range("M4").Value = -128
$ws3.Range("H43").Value = 0
$ws3.Range("J43").Value = 0
$ws3.Range("L43").Value = 0
$ws3.Range("N43").ClearContents()

$ws4 = $wb.Worksheets.Item("CRP")
$ws4.Range("H22").Value = 909.8
$ws4.Range("I22").Value = 866.3333
$ws4.Range("K22").Value = 866.3333
$ws4.Range("M22").Value = -516.3333
$ws4.Range("H31").Value = 3079.1184
$ws4.Range("I31").Value = 2067.9375
$ws4.Range("K31").Value = 2067.9375
$ws4.Range("M31").Value = -1772.9375
$ws4.Range("H34").Value = 3079.1184
$ws4.Range("I34").Value = 2067.9375
$ws4.Range("K34").Value = 2067.9375
$ws4.Range("M34").Value = -1865.9375

$ws5 = $wb.Worksheets.Item("CUL")
$ws5.Range("H56").Value = 2750
$ws5.Range("I56").Value = 2750
$ws5.Range("K56").Value = 2750
$ws5.Range("M56").Value = -2220
$ws5.Range("H69").Value = 1924.3334
$ws5.Range("J69").Value = 2000
$ws5.Range("L69").Value = 6000
$ws5.Range("N69").Value = -7622
$ws5.Range("H72").Value = 1924.3334
$ws5.Range("J72").Value = 2000
$ws5.Range("L72").Value = 18000
$ws5.Range("N72").Value = -26112
$ws5.Range("H81").Value = 6762.6665
$ws5.Range("I81").Value = 5000
$ws5.Range("K81").Value = 15000
$ws5.Range("M81").Value = -13877
$ws5.Range("H84").Value = 6762.6665
$ws5.Range("I84").Value = 5000
$ws5.Range("K84").Value = 45000
$ws5.Range("M84").Value = -39384
$ws5.Range("H99").Value = 5199.9
$ws5.Range("I99").Value = 1999.5
$ws5.Range("K99").Value = 5998.5
$ws5.Range("M99").Value = -3752.5
$ws5.Range("H102").Value = 19362.5
$ws5.Range("I102").Value = 14150
$ws5.Range("K102").Value = 42450
$ws5.Range("M102").Value = -40016
$ws5.Range("H105").Value = 12582.667
$ws5.Range("I105").Value = 4000
$ws5.Range("J105").Value = 14299.2
$ws5.Range("K105").Value = 12000
$ws5.Range("L105").Value = 42897.60000000001
$ws5.Range("M105").Value = -9379
$ws5.Range("N105").Value = -48139.60000000001
$ws5.Range("H132").Value = 2756.2534
$ws5.Range("I132").Value = 1774.8334
$ws5.Range("J132").Value = 2955.8645
$ws5.Range("K132").Value = 15973.5006
$ws5.Range("L132").Value = 26602.7805
$ws5.Range("M132").Value = -13443.5006
$ws5.Range("N132").Value = -31662.7805

$ws6 = $wb.Worksheets.Item("GSM")
$ws6.Range("H80").Value = 3524.4
$ws6.Range("J80").Value = 4733.154
$ws6.Range("L80").Value = 4733.154
$ws6.Range("N80").Value = -6729.154
$ws6.Range("H83").Value = 3524.4
$ws6.Range("J83").Value = 4733.154
$ws6.Range("L83").Value = 23665.77
$ws6.Range("N83").Value = -33649.77

$ws7 = $wb.Worksheets.Item("LTW")
$ws7.Range("H7").Value = 3234.76
$ws7.Range("I7").Value = 2957.1177
$ws7.Range("J7").Value = 3824.75
$ws7.Range("K7").Value = 2957.1177
$ws7.Range("L7").Value = 3824.75
$ws7.Range("M7").Value = -2845.1177
$ws7.Range("N7").Value = -4048.75
$ws7.Range("H46").Value = 1921.7
$ws7.Range("J46").Value = 2523.5715
$ws7.Range("L46").Value = 2523.5715
$ws7.Range("N46").Value = -2899.5715
$ws7.Range("H61").Value = 3525.4167
$ws7.Range("I61").Value = 4917.1665
$ws7.Range("J61").Value = 2133.6667
$ws7.Range("K61").Value = 4917.1665
$ws7.Range("L61").Value = 2133.6667
$ws7.Range("M61").Value = -4715.1665
$ws7.Range("N61").Value = -2537.6667
$ws7.Range("H82").Value = 750.5
$ws7.Range("I82").Value = 750.5
$ws7.Range("J82").Value = 0
$ws7.Range("K82").Value = 750.5
$ws7.Range("L82").Value = 0
$ws7.Range("M82").Value = -389.5
$ws7.Range("N82").ClearContents()
$ws7.Range("H85").Value = 750.5
$ws7.Range("I85").Value = 750.5
$ws7.Range("J85").Value = 0
$ws7.Range("K85").Value = 750.5
$ws7.Range("L85").Value = 0
$ws7.Range("M85").Value = 497.5
$ws7.Range("N85").ClearContents()
$ws7.Range("H113").Value = 3525.4167
$ws7.Range("I113").Value = 4917.1665
$ws7.Range("J113").Value = 2133.6667
$ws7.Range("K113").Value = 4917.1665
$ws7.Range("L113").Value = 2133.6667
$ws7.Range("M113").Value = -2747.1665
$ws7.Range("N113").Value = -6473.6667
$ws7.Range("H126").Value = 3234.76
$ws7.Range("I126").Value = 2957.1177
$ws7.Range("J126").Value = 3824.75
$ws7.Range("K126").Value = 8871.3531
$ws7.Range("L126").Value = 11474.25
$ws7.Range("M126").Value = -6401.3531
$ws7.Range("N126").Value = -16414.25
$ws7.Range("H132").Value = 31253230
$ws7.Range("I132").Value = 55558516
$ws7.Range("K132").Value = 166675548
$ws7.Range("M132").Value = -166673018

$ws8 = $wb.Worksheets.Item("WVR")
$ws8.Range("H58").Value = 19999.334
$ws8.Range("I58").Value = 19999
$ws8.Range("K58").Value = 19999
$ws8.Range("M58").Value = -19691
$ws8.Range("H107").Value = 1119.5454
$ws8.Range("I107").Value = 1085.6666
$ws8.Range("J107").Value = 1160.2
$ws8.Range("K107").Value = 3256.9998
$ws8.Range("L107").Value = 3480.6
$ws8.Range("M107").Value = -1336.9998
$ws8.Range("N107").Value = -7320.6
$ws8.Range("H132").Value = 8848.611000000001
$ws8.Range("I132").Value = 10935.091
$ws8.Range("J132").Value = 5569.857
$ws8.Range("K132").Value = 32805.273
$ws8.Range("L132").Value = 16709.571
$ws8.Range("M132").Value = -30275.273
$ws8.Range("N132").Value = -21769.571
